# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-17 Monday", 2) | Out-Null

# Update the division-problem answers in the single table.
# Each content row in the table (1-based Word row numbers) holds 5 answers;
# addressing cells directly by (row, column) avoids any ambiguity from
# duplicate cell text (e.g. "61÷4=15, 1" appears twice with different
# replacements).
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Cell.Range includes the trailing end-of-cell mark; trim it off so we
    # only replace the run's visible text (preserving its formatting).
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = $newText
}

Set-CellText $t 1 1 "88÷6=14, 4"
Set-CellText $t 1 2 "36÷8=4, 4"
Set-CellText $t 1 3 "77÷4=19, 1"
Set-CellText $t 1 4 "18÷2=9, 0"
Set-CellText $t 1 5 "84÷5=16, 4"

Set-CellText $t 5 1 "90÷2=45, 0"
Set-CellText $t 5 2 "24÷7=3, 3"
Set-CellText $t 5 3 "41÷8=5, 1"
Set-CellText $t 5 4 "81÷4=20, 1"
Set-CellText $t 5 5 "92÷2=46, 0"

Set-CellText $t 9 1 "24÷3=8, 0"
Set-CellText $t 9 2 "43÷9=4, 7"
Set-CellText $t 9 3 "86÷7=12, 2"
Set-CellText $t 9 4 "43÷8=5, 3"
Set-CellText $t 9 5 "95÷3=31, 2"

Set-CellText $t 13 1 "76÷3=25, 1"
Set-CellText $t 13 2 "19÷6=3, 1"
Set-CellText $t 13 3 "45÷4=11, 1"
Set-CellText $t 13 4 "48÷3=16, 0"
Set-CellText $t 13 5 "28÷6=4, 4"

Set-CellText $t 17 1 "59÷3=19, 2"
Set-CellText $t 17 2 "85÷9=9, 4"
Set-CellText $t 17 3 "73÷2=36, 1"
Set-CellText $t 17 4 "11÷8=1, 3"
Set-CellText $t 17 5 "51÷5=10, 1"
